$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the existing content (rows 2 and 3),
# pushing everything else down by 2 rows.
$ws.Rows("2:3").Insert()

# Add the new title in B2, matching the bold/number-format style used
# elsewhere in the sheet (same style as the SUM cell).
$ws.Range("B2").Value = "Project Euler 2: Even Fibonacci Numbers"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").NumberFormat = "#,##0_ ;[Red]\-#,##0\ "

# Restore the originally selected cell.
$ws.Range("B3").Select()
